# BP-1588 category all apis
# Adds a "category" column to the "Курсы" sheet (new column E), populated
# with "Профориентация" for every course row, and makes "Курсы" the
# active sheet/tab again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Курсы")

# Insert a new column before the current "type" column (E), shifting
# type/options one column to the right (E->F, F->G) and copying
# formatting the way Excel's own Insert does.
$ws.Columns("E").Insert() | Out-Null

# Header + data for the new "category" column.
$ws.Range("E1").Value = "category"
$ws.Range("E2").Value = "Профориентация"
$ws.Range("E3").Value = "Профориентация"
$ws.Range("E4").Value = "Профориентация"
$ws.Range("E5").Value = "Профориентация"

# Match the new column's width to the target layout.
$ws.Columns("E").ColumnWidth = 18

# Page setup (A4 / portrait) shows up on this sheet in the saved file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore "Курсы" as the active/selected sheet + cell so the workbook
# reopens on it instead of "Опции".
$ws.Activate() | Out-Null
$ws.Range("E1:E2").Select() | Out-Null
